# Append a new logged experiment run (row 18) at the bottom of the
# experiments table, mirroring the existing rows' layout:
# Function | Iterations | Positive | Negative | Error | Duration | Iters per Second
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 18
$ws.Cells.Item($row, 1).Value = "x**2"
$ws.Cells.Item($row, 2).Value = 20000
$ws.Cells.Item($row, 3).Value = 1
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0.0004686458851210773
$ws.Cells.Item($row, 6).Value = 22.45450019836426
$ws.Cells.Item($row, 7).Value = 890.6900542572281
